$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 7497256
$ws.Range("C4").Value = 2585
$ws.Range("D4").Value = 4737369
$ws.Range("E4").Value = 2547193
$ws.Range("G4").Value = 34
$ws.Range("H4").Value = 212694

# Row 16: Iran
$ws.Range("A16").Value = "Iran"
$ws.Range("B16").Value = 464596
$ws.Range("C16").Value = 3552
$ws.Range("D16").Value = 385264
$ws.Range("E16").Value = 52765
$ws.Range("G16").Value = 187
$ws.Range("H16").Value = 26567

# Row 25: Alemania
$ws.Range("A25").Value = "Alemania"
$ws.Range("B25").Value = 295943
$ws.Range("C25").Value = 413
$ws.Range("D25").Value = 259500
$ws.Range("E25").Value = 26857

# Row 32: Rumania
$ws.Range("A32").Value = "Rumania"
$ws.Range("B32").Value = 132001
$ws.Range("C32").Value = 2343
$ws.Range("D32").Value = 105582
$ws.Range("E32").Value = 21504
$ws.Range("G32").Value = 53
$ws.Range("H32").Value = 4915

# Row 57: Barein
$ws.Range("A57").Value = "Barein"
$ws.Range("E57").Value = 5569
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = 255

# Row 61: Suiza
$ws.Range("A61").Value = "Suiza"
$ws.Range("B61").Value = 54384
$ws.Range("C61").Value = 552
$ws.Range("D61").Value = 45300
$ws.Range("E61").Value = 7009
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 2075

# Row 62: Moldavia
$ws.Range("A62").Value = "Moldavia"
$ws.Range("B62").Value = 54064
$ws.Range("D62").Value = 39499
$ws.Range("E62").Value = 13229
$ws.Range("H62").Value = 1336

# Row 73: Kenia
$ws.Range("A73").Value = "Kenia"
$ws.Range("D73").Value = 25023
$ws.Range("E73").Value = 12972
$ws.Range("H73").Value = 718

# Row 90: Madagascar
$ws.Range("A90").Value = "Madagascar"
$ws.Range("B90").Value = 16493
$ws.Range("C90").Value = 39
$ws.Range("D90").Value = 15467
$ws.Range("E90").Value = 794

# Row 91: Senegal
$ws.Range("A91").Value = "Senegal"
$ws.Range("B91").Value = 15051
$ws.Range("C91").Value = 32
$ws.Range("D91").Value = 12694
$ws.Range("E91").Value = 2045
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 312

# Row 98: Malasia
$ws.Range("A98").Value = "Malasia"
$ws.Range("B98").Value = 11771
$ws.Range("C98").Value = 287
$ws.Range("D98").Value = 10095
$ws.Range("E98").Value = 1540
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 136

# Row 99: Eslovaquia
$ws.Range("A99").Value = "Eslovaquia"
$ws.Range("B99").Value = 11617
$ws.Range("C99").Value = 679
$ws.Range("D99").Value = 4756
$ws.Range("E99").Value = 6807
$ws.Range("G99").Value = 6
$ws.Range("H99").Value = 54

# Row 142: Estonia
$ws.Range("A142").Value = "Estonia"
$ws.Range("E142").Value = 766
$ws.Range("G142").Value = 1
$ws.Range("H142").Value = 66

# Row 145: Malta
$ws.Range("A145").Value = "Malta"
$ws.Range("B145").Value = 3139
$ws.Range("C145").Value = 44
$ws.Range("D145").Value = 2668
$ws.Range("E145").Value = 434
$ws.Range("G145").Value = 2
$ws.Range("H145").Value = 37

# Row 146: Mali
$ws.Range("A146").Value = "Mali"
$ws.Range("B146").Value = 3131
$ws.Range("D146").Value = 2460
$ws.Range("E146").Value = 540
$ws.Range("H146").Value = 131

# Row 180: Gibraltar
$ws.Range("A180").Value = "Gibraltar"
$ws.Range("B180").Value = 416
$ws.Range("C180").Value = 6
$ws.Range("D180").Value = 351
$ws.Range("E180").Value = 65

# Row 196: Liechtenstein
$ws.Range("A196").Value = "Liechtenstein"
$ws.Range("B196").Value = 120
$ws.Range("C196").Value = 1
$ws.Range("D196").Value = 116
$ws.Range("E196").Value = 3

# Row 207: Santa Lucia
$ws.Range("A207").Value = "Santa Lucia"

# Row 208: Nueva Caledonia
$ws.Range("A208").Value = "Nueva Caledonia"

# Update timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 2 de Octubre de 2020 a las 13:01"
